$d = $word.ActiveDocument
$enDash = [char]0x2013

# 1. Remove the old "_GoBack" bookmark (it sits in the "Menu :" bullet paragraph).
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 2. In the "Highscore" bullet, swap the order of the "score" and "gold" tokens:
#    "nama player - gold [enDash] score - level ..."
#      becomes
#    "nama player [enDash] score - gold - level ..."
$scope = $d.Content
$scope.Find.Execute("nama player")
$afterNamaPlayer = $scope.End
$tail = $d.Range($afterNamaPlayer, $d.Content.End)

$oldText = " - gold " + $enDash + " score - "
$newText = " " + $enDash + " score - gold - "
$tail.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# 3. Re-insert the "_GoBack" bookmark right after the (now relocated) "gold" token.
$scope2 = $d.Content
$scope2.Find.Execute("nama player")
$afterNamaPlayer2 = $scope2.End
$goldRange = $d.Range($afterNamaPlayer2, $d.Content.End)
$goldRange.Find.Execute("gold")
$bookmarkPos = $goldRange.End
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
